# Update the email address stored on the "Email" sheet (cell B2) and make
# that sheet the active one, matching the author's change "Email address in
# config file updated".

$wb = $excel.ActiveWorkbook

$wsEmail = $wb.Worksheets.Item("Email")

# Turn B2 into a mailto hyperlink pointing at the new address first (adding
# the hyperlink before touching the cell's value/style keeps the existing
# built-in "Hyperlink" cell style (index 1) instead of Excel registering a
# duplicate one for the cell).
$wsEmail.Hyperlinks.Add($wsEmail.Range("B2"), "mailto:delia.panca@fwfcompany.com") | Out-Null

# Set the visible text and apply the standard Hyperlink style.
$wsEmail.Range("B2").Value = "delia.panca@fwfcompany.com"
$wsEmail.Range("B2").Style = "Hyperlink"

# Make the Email sheet the active/selected tab, with B6 as the selected cell.
$wsEmail.Activate() | Out-Null
$wsEmail.Range("B6").Select() | Out-Null
